# Update experiment values in row 2 (ti_MM results) with improved CI values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.556060110305386
$ws.Range("E2").Value = 2.197559776593686
$ws.Range("F2").Value = 2.835601846558793
$ws.Range("G2").Value = 3.433433688526137
$ws.Range("H2").Value = 3.972124914913886
$ws.Range("I2").Value = 4.443834622851361
$ws.Range("J2").Value = 4.847085175494485
$ws.Range("K2").Value = 5.183430842096524
$ws.Range("L2").Value = 5.455506157444306
$ws.Range("M2").Value = 5.658309842224586
$ws.Range("N2").Value = 5.795671990345435
$ws.Range("O2").Value = 5.869153917792468
$ws.Range("P2").Value = 5.877817653592294
$ws.Range("Q2").Value = 5.833945766326782
$ws.Range("R2").Value = 5.760697970840522
$ws.Range("S2").Value = 5.673253170923223
$ws.Range("T2").Value = 5.581356806311221
$ws.Range("U2").Value = 5.491067691967056
$ws.Range("V2").Value = 5.405958968213211
$ws.Range("W2").Value = 5.327945106338505
$ws.Range("X2").Value = 5.257851493563179
$ws.Range("Y2").Value = 5.195805961433192
$ws.Range("Z2").Value = 5.141506650078503
$ws.Range("AA2").Value = 5.094403744668927
$ws.Range("AB2").Value = 5.053821167629343
$ws.Range("AC2").Value = 5.019036464908085
$ws.Range("AD2").Value = 4.989331701742379
$ws.Range("AE2").Value = 4.964024400924715
$ws.Range("AF2").Value = 4.946251086557831

$wb.Save()
